$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.6260102880529633
$ws.Range("D3").Value = 0.3541428331643913
$ws.Range("D4").Value = 0.679786453153583
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.4986321910239432
$ws.Range("D6").Value = 0.4223807412409279
$ws.Range("D7").Value = 0.7062903966157175
$ws.Range("D8").Value = 0.5459168673440993
$ws.Range("D9").Value = 0.3137614743579401
$ws.Range("D10").Value = 0.4780994288015241
$ws.Range("D11").Value = 0.3508560201085118
$ws.Range("D12").Value = 0.4195763179118504
$ws.Range("D13").Value = 0.7143835820507201
$ws.Range("D14").Value = 0.5371074731704304
$ws.Range("D15").Value = 0.6091842103336705
$ws.Range("D16").Value = 0.4142718674571033
$ws.Range("D17").Value = 0.678384174770905
$ws.Range("D18").Value = 0.6043942648274618
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 0.5763493757467612
$ws.Range("D20").Value = 0.4405107685376133
$ws.Range("D21").Value = 0.4022043166551263
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0.498489480808465
$ws.Range("D23").Value = 0.6333580243952991
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0.4346929416979716
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 0.5957396472105976
$ws.Range("D26").Value = 0.7167864482887295
$ws.Range("D27").Value = 0.3028307437312055
$ws.Range("D28").Value = 0.5640891306221512
$ws.Range("D29").Value = 0.3562952915245866
$ws.Range("D30").Value = 0.5516865963826384
$ws.Range("D31").Value = 0.5798070832243011
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 0.4999938270853297
$ws.Range("D33").Value = 0.4892916420992158
$ws.Range("D34").Value = 0.7216873612678408
$ws.Range("D35").Value = 0.5737616424322874
$ws.Range("D36").Value = 0.4711000822257232
$ws.Range("D37").Value = 0.6526260348467665
$ws.Range("D38").Value = 0.5912504099833619
$ws.Range("D39").Value = 0.6618911604902618
$ws.Range("D40").Value = 0.6990647244019231
$ws.Range("C41").Value = 1
$ws.Range("D41").Value = 0.6084400649470929
$ws.Range("D42").Value = 0.5090362921898108
$ws.Range("D43").Value = 0.539152677276871
$ws.Range("C44").Value = 1
$ws.Range("D44").Value = 0.5674482865947662
$ws.Range("D45").Value = 0.4359300577359477
$ws.Range("D46").Value = 0.6442761949792856
$ws.Range("D47").Value = 0.536265694211266
$ws.Range("C48").Value = 1
$ws.Range("D48").Value = 0.5340428357068283
$ws.Range("D49").Value = 0.5129239904461157
$ws.Range("D50").Value = 0.4704472608223459
$ws.Range("D51").Value = 0.6248170270389944
$ws.Range("D52").Value = 0.3137539086539305
$ws.Range("D53").Value = 0.6189485826895802
$ws.Range("D54").Value = 0.633963277300522
$ws.Range("D55").Value = 0.3312994607613299
$ws.Range("D56").Value = 0.4956921585003743
$ws.Range("D57").Value = 0.456294799912472
$ws.Range("D58").Value = 0.6424620692710761
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 0.4066792766240554
$ws.Range("D60").Value = 0.6464608402776206
$ws.Range("D61").Value = 0.5302003615961089
$ws.Range("D62").Value = 0.6956237122153737
$ws.Range("D63").Value = 0.4952512695307403
$ws.Range("D64").Value = 0.3863297190865318
$ws.Range("D65").Value = 0.3570019636339106
$ws.Range("D66").Value = 0.5545425254074483
$ws.Range("C67").Value = 0
$ws.Range("D67").Value = 0.4980158376756523
$ws.Range("C68").Value = 1
$ws.Range("D68").Value = 0.5007556900739473
$ws.Range("D69").Value = 0.5862153277618886
$ws.Range("D70").Value = 0.3356683018530139
$ws.Range("C71").Value = 1
$ws.Range("D71").Value = 0.5013643205764443
$ws.Range("C72").Value = 1
$ws.Range("D72").Value = 0.5002903240965598
$ws.Range("D73").Value = 0.5432051249583271
$ws.Range("D74").Value = 0.373331320770628
$ws.Range("D75").Value = 0.3877931438908669
$ws.Range("D76").Value = 0.3585776163144284
$ws.Range("D77").Value = 0.5061076558188958
$ws.Range("D78").Value = 0.6171779257778169
$ws.Range("D79").Value = 0.4410681770014703
